# Slovenia Prva Liga - base update (19-04-2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 9 and 10: the two match records were swapped (same columns A/C/D/E,
# everything else exchanged between the rows).
# ---------------------------------------------------------------------------
$ws.Cells.Item(9, 2).Value  = 6814328
$ws.Cells.Item(9, 6).Value  = "NK Domzale"
$ws.Cells.Item(9, 7).Value  = "NK Bravo"
$ws.Cells.Item(9, 8).Value  = 1
$ws.Cells.Item(9, 9).Value  = 1
$ws.Cells.Item(9, 10).Value = "D"
$ws.Cells.Item(9, 11).Value = 2.35
$ws.Cells.Item(9, 12).Value = 3.1
$ws.Cells.Item(9, 13).Value = 2.9
$ws.Cells.Item(9, 14).Value = 2.15
$ws.Cells.Item(9, 15).Value = 3.1
$ws.Cells.Item(9, 16).Value = 3.3
$ws.Cells.Item(9, 17).Value = -0.25
$ws.Cells.Item(9, 18).Value = 1.925
$ws.Cells.Item(9, 19).Value = 1.875
$ws.Cells.Item(9, 20).Value = 2.25
$ws.Cells.Item(9, 21).Value = 1.95
$ws.Cells.Item(9, 22).Value = 1.85
$ws.Cells.Item(9, 23).Value = -1
$ws.Cells.Item(9, 24).Value = 2.1
$ws.Cells.Item(9, 25).Value = -1
$ws.Cells.Item(9, 26).Value = -0.5
$ws.Cells.Item(9, 27).Value = 0.4375
$ws.Cells.Item(9, 28).Value = -0.5
$ws.Cells.Item(9, 29).Value = 0.425

$ws.Cells.Item(10, 2).Value  = 6814330
$ws.Cells.Item(10, 6).Value  = "NK Maribor"
$ws.Cells.Item(10, 7).Value  = "NK Aluminij"
$ws.Cells.Item(10, 8).Value  = 1
$ws.Cells.Item(10, 9).Value  = 0
$ws.Cells.Item(10, 10).Value = "H"
$ws.Cells.Item(10, 11).Value = 1.363
$ws.Cells.Item(10, 12).Value = 4.5
$ws.Cells.Item(10, 13).Value = 7
$ws.Cells.Item(10, 14).Value = 1.4
$ws.Cells.Item(10, 15).Value = 4.5
$ws.Cells.Item(10, 16).Value = 7
$ws.Cells.Item(10, 17).Value = -1.25
$ws.Cells.Item(10, 18).Value = 1.85
$ws.Cells.Item(10, 19).Value = 1.95
$ws.Cells.Item(10, 20).Value = 2.75
$ws.Cells.Item(10, 21).Value = 1.8
$ws.Cells.Item(10, 22).Value = 2
$ws.Cells.Item(10, 23).Value = 0.3999999999999999
$ws.Cells.Item(10, 24).Value = -1
$ws.Cells.Item(10, 25).Value = -1
$ws.Cells.Item(10, 26).Value = -0.5
$ws.Cells.Item(10, 27).Value = 0.475
$ws.Cells.Item(10, 28).Value = -1
$ws.Cells.Item(10, 29).Value = 1

# ---------------------------------------------------------------------------
# Row 152: this fixture record was entirely replaced with a new match.
# ---------------------------------------------------------------------------
$ws.Cells.Item(152, 2).Value  = 6814424
$ws.Cells.Item(152, 5).Value  = 45399.5625
$ws.Cells.Item(152, 6).Value  = "Olimpija Ljubljana"
$ws.Cells.Item(152, 7).Value  = "NK Maribor"
$ws.Cells.Item(152, 8).Value  = 1
$ws.Cells.Item(152, 9).Value  = 2
$ws.Cells.Item(152, 10).Value = "A"
$ws.Cells.Item(152, 11).Value = 2.375
$ws.Cells.Item(152, 12).Value = 3.4
$ws.Cells.Item(152, 13).Value = 2.65
$ws.Cells.Item(152, 14).Value = 2.1
$ws.Cells.Item(152, 15).Value = 3.6
$ws.Cells.Item(152, 16).Value = 3
$ws.Cells.Item(152, 17).Value = -0.25
$ws.Cells.Item(152, 18).Value = 1.85
$ws.Cells.Item(152, 19).Value = 1.95
$ws.Cells.Item(152, 20).Value = 3
$ws.Cells.Item(152, 21).Value = 2
$ws.Cells.Item(152, 22).Value = 1.8
$ws.Cells.Item(152, 23).Value = -1
$ws.Cells.Item(152, 24).Value = -1
$ws.Cells.Item(152, 25).Value = 2
$ws.Cells.Item(152, 26).Value = -1
$ws.Cells.Item(152, 27).Value = 0.95
$ws.Cells.Item(152, 28).Value = 0
$ws.Cells.Item(152, 29).Value = -0

# ---------------------------------------------------------------------------
# Row 153: closing-odds update (U/V only).
# ---------------------------------------------------------------------------
$ws.Cells.Item(153, 21).Value = 1.775
$ws.Cells.Item(153, 22).Value = 2.025

# ---------------------------------------------------------------------------
# Row 154: closing-odds update (N/O/P and R/S/U/V).
# ---------------------------------------------------------------------------
$ws.Cells.Item(154, 14).Value = 1.333
$ws.Cells.Item(154, 15).Value = 4.75
$ws.Cells.Item(154, 16).Value = 7.5
$ws.Cells.Item(154, 18).Value = 1.975
$ws.Cells.Item(154, 19).Value = 1.825
$ws.Cells.Item(154, 21).Value = 2
$ws.Cells.Item(154, 22).Value = 1.8

# ---------------------------------------------------------------------------
# Row 155: brand-new fixture appended at the bottom of the table. Clone the
# number formatting/borders from row 154 (id/date columns carry styles) then
# fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A154").Copy() | Out-Null
$ws.Range("A155").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("E154").Copy() | Out-Null
$ws.Range("E155").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Cells.Item(155, 1).Value  = 153
$ws.Cells.Item(155, 2).Value  = 6994888
$ws.Cells.Item(155, 3).Value  = "Slovenia Prva Liga"
$ws.Cells.Item(155, 4).Value  = "Slovenia Prva Liga"
$ws.Cells.Item(155, 5).Value  = 45403.63541666666
$ws.Cells.Item(155, 6).Value  = "FC Koper"
$ws.Cells.Item(155, 7).Value  = "NK Bravo"
$ws.Cells.Item(155, 11).Value = 2.2
$ws.Cells.Item(155, 12).Value = 3.3
$ws.Cells.Item(155, 13).Value = 3
$ws.Cells.Item(155, 14).Value = 2.15
$ws.Cells.Item(155, 15).Value = 3.3
$ws.Cells.Item(155, 16).Value = 3.1
$ws.Cells.Item(155, 17).Value = -0.25
$ws.Cells.Item(155, 18).Value = 1.9
$ws.Cells.Item(155, 19).Value = 1.9
$ws.Cells.Item(155, 20).Value = 2.25
$ws.Cells.Item(155, 21).Value = 1.9
$ws.Cells.Item(155, 22).Value = 1.9
$ws.Cells.Item(155, 23).Value = 0
$ws.Cells.Item(155, 24).Value = 0
$ws.Cells.Item(155, 25).Value = 0
$ws.Cells.Item(155, 26).Value = 0
$ws.Cells.Item(155, 27).Value = 0
